$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 873, shifting existing rows 873+ down by one.
$ws.Rows("873:873").Insert()

# Populate the newly inserted row with the new data point.
# Force column A to remain plain text (matches the rest of the date column,
# which is stored as text, not as a date serial number) and keep it
# unstyled, same as the other data rows.
$ws.Range("A873").NumberFormat = "@"
$ws.Range("A873").Value = "2026/02/24"
$ws.Range("A873").ClearFormats()
$ws.Range("B873").Value = "火"
$ws.Range("C873").Value = 12
$ws.Range("D873").Value = 201
